$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.947.08'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.664.48'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.68'
$ws.Range("E5").Value = '  +5.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3632'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.79'
$ws.Range("E8").Value = '  +2.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3274'
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.136'
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07101'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9996'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.066'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.60'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '1.661.89'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.627'
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001050'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06614'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9988'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '79.23'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.931'
$ws.Range("E21").Value = '  -1.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.83'
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.62'
$ws.Range("E23").Value = '  +3.27%  '
$ws.Range("D24").Value = '24.899.33'
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.459'
$ws.Range("E25").Value = '  +1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.442'
$ws.Range("E26").Value = '  -3.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.00'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("D29").Value = '1.842.60'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.70'
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.185'
$ws.Range("E31").Value = '  +7.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.092'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.757'
$ws.Range("E33").Value = '  -5.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08479'
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.657'
$ws.Range("E35").Value = '  -4.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.29'
$ws.Range("E36").Value = '  -2.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.285'
$ws.Range("E37").Value = '  +5.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.192'
$ws.Range("E38").Value = '  -0.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06152'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02276'
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.321'
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2079'
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9995'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5949'
$ws.Range("E44").Value = '  -1.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.63'
$ws.Range("E45").Value = '  +4.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.836'
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5642'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.35'
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.954'
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06994'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.192'
$ws.Range("E51").Value = '  +2.02%  '
